# Update annotations for Ruilin
# - Row 75, column B (politeness_score): change stored type from text "3" to numeric 3
# - Insert a new row 76 with a fresh annotation record

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 75: make B75 a real number (was text "3") -------------------------
$ws.Range("B75").Value = 3

# --- Row 76: new annotation row ---------------------------------------------
$ws.Range("A76").Value = "Ruilin"

# B76 must stay a *text* "3" (matches style of the original B75 before the
# type fix above) - force text entry via NumberFormat, then strip the format
# back to Normal so no stray style sticks around on the cell.
$ws.Range("B76").NumberFormat = "@"
$ws.Range("B76").Value = "3"
$ws.Range("B76").Style = "Normal"

$ws.Range("C76").Value = "无"
$ws.Range("D76").Value = "DIS"
$ws.Range("E76").Value = "OTH"
$ws.Range("F76").Value = "85844681-e6c1-4472-a9f5-69a1244b25a4"
$ws.Range("G76").Value = "SktLlGbRZ_annotated.xlsx"
$ws.Range("H76").Value = "Also, it would be good to extend the figure with the second cycle loss."
